$wb = $excel.ActiveWorkbook

# --- Rebuild "Sensor Data" sheet layout ---
$ws = $wb.Worksheets.Item("Sensor Data")

# Clear the old content (A1:F11) before rebuilding with the new column layout.
$ws.Range("A1:F11").Clear()

# Header row: base_addr | variable_name | addr_offset | type | init_value | hmi_tag
$ws.Range("A1").Value = "base_addr"
$ws.Range("B1").Value = "variable_name"
$ws.Range("C1").Value = "addr_offset"
$ws.Range("D1").Value = "type"
$ws.Range("E1").Value = "init_value"
$ws.Range("F1").Value = "hmi_tag"
$ws.Range("B1:F1").HorizontalAlignment = -4108

$ws.Range("B2").Value = "value"
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "WORD"
$ws.Range("E2").Value = 0

$ws.Range("B3").Value = "state"
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "WORD"
$ws.Range("E3").Value = 0

$ws.Range("B4").Value = "err_u"
$ws.Range("C4").Value = 1
$ws.Range("D4").Value = "WORD"
$ws.Range("E4").Value = 0

$ws.Range("B5").Value = "warn_u"
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = "WORD"
$ws.Range("E5").Value = 0

$ws.Range("B6").Value = "warn_l"
$ws.Range("C6").Value = 1
$ws.Range("D6").Value = "WORD"
$ws.Range("E6").Value = 0

$ws.Range("B7").Value = "err_l"
$ws.Range("C7").Value = 1
$ws.Range("D7").Value = "WORD"
$ws.Range("E7").Value = 0

# addr_offset column is right-aligned (new style picked up by row 2-8)
$ws.Range("C2:C8").HorizontalAlignment = -4152

# F3 keeps the centered style carried over from the old table footer
$ws.Range("F3").HorizontalAlignment = -4108

# trailing helper row (blank, addr_offset style carried down)
$ws.Range("C8").HorizontalAlignment = -4152

# Column widths to match the new, narrower data
$ws.Columns.Item(1).ColumnWidth = 8.83203125
$ws.Columns.Item(2).ColumnWidth = 12.1640625
$ws.Columns.Item(3).ColumnWidth = 9.83203125
$ws.Columns.Item(4).ColumnWidth = 6
$ws.Columns.Item(5).ColumnWidth = 8.5
$ws.Columns.Item(6).ColumnWidth = 7.1640625

# --- Selection / active-tab bookkeeping ---
$wsConstants = $wb.Worksheets.Item("Constants")
$wsConstants.Range("A1:F1").Select()

$ws.Activate()
$ws.Range("G5").Select()
